# IEEEST now supports a remote bus specified by `busr`.
#
# 1. Reorder the sheet tabs so "IEEEST" appears before "Toggler"
#    (previously Toggler, IEEEST -> now IEEEST, Toggler).
# 2. Insert the new `busr` column into the IEEEST sheet, right after
#    MODE (and before A1), and populate its header + sample value.
# 3. Update the sample MODE value and fill in the previously-missing
#    uid value on the IEEEST data row.

$wb = $excel.ActiveWorkbook

$ieeest  = $wb.Worksheets.Item("IEEEST")
$toggler = $wb.Worksheets.Item("Toggler")
$ieeest.Move($toggler)

$ws = $wb.Worksheets.Item("IEEEST")

# Insert the busr column (currently column G holds "A1").
$ws.Columns.Item(7).Insert()
$ws.Cells.Item(1, 7).Value = "busr"
$ws.Cells.Item(2, 7).Value = 2

# MODE (column F) sample changes from 1 -> 5.
$ws.Cells.Item(2, 6).Value = 5

# uid (column A) sample value was missing; now 0.
$ws.Cells.Item(2, 1).Value = 0

# Keep IEEEST as the active/selected sheet, matching the source edit.
$ws.Activate()
